# "Add cantrals by cantons"
# Turn the two-row (partial) header of Sheet1 into a single, fully-labelled
# header row and remove the now-redundant second header row, shifting all
# the plant/canton data rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Drop the old sub-header row (Hiver/Eté/Année + (m3/s)/(MW)/(GWh)); this
#    shifts every data row up by one (old row 3 -> new row 2, etc.)
$ws.Rows(2).Delete()

# 2) Write the new, single header row.
$ws.Range("A1").Value2 = "idx"
$ws.Range("B1").Value2 = "idx2"
$ws.Range("C1").Value2 = "Name"
$ws.Range("D1").Value2 = "Date Start"
$ws.Range("E1").Value2 = "Date End"
$ws.Range("F1").Value2 = "(m3/s)"
$ws.Range("G1").Value2 = "(MW1)"
$ws.Range("H1").Value2 = "(MW2)"
$ws.Range("I1").Value2 = "(GWh) Winter"
$ws.Range("J1").Value2 = "(GWh) Summer"
$ws.Range("K1").Value2 = "(GWh) Year"

# 3) idx/idx2/Name/Date Start/Date End header cells use the plain default
#    style (no special formatting).
$ws.Range("A1:E1").Style = "Normal"

# 4) (m3/s)/(MW1)/(MW2)/(GWh)... header cells use the small (9pt Arial)
#    font used throughout the rest of the sheet, with no special number
#    format. Add a throwaway named style to get that exact combination,
#    apply it, then remove the named style again (the cell keeps the
#    resulting cell format).
$headerStyle = $wb.Styles.Add("HeaderStyle")
$headerStyle.Font.Name = "Arial"
$headerStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "HeaderStyle"
$wb.Styles.Item("HeaderStyle").Delete()

# 5) Match the selection left behind by the edit (first data row selected).
$ws.Range("A2:K2").Select()

Write-Output "done"
